$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.3371691897672
$ws.Cells.Item(2, 3).Value = 9.90525935230265
$ws.Cells.Item(2, 4).Value = 9.585710157139934
$ws.Cells.Item(2, 6).Value = 33.88906282279513
$ws.Cells.Item(2, 7).Value = 35.45348869832786
$ws.Cells.Item(2, 8).Value = 16.1438656118313
$ws.Cells.Item(2, 9).Value = 26.02771544719502
$ws.Cells.Item(2, 10).Value = 10.94858464364165
$ws.Cells.Item(2, 11).Value = 8.36243219960976
$ws.Cells.Item(2, 12).Value = 11.65678769520694
$ws.Cells.Item(2, 13).Value = 14.63566406604197
$ws.Cells.Item(2, 15).Value = 25.42685309591087

$ws.Cells.Item(3, 2).Value = 11.08306694070348
$ws.Cells.Item(3, 3).Value = 9.896947129081033
$ws.Cells.Item(3, 4).Value = 9.575887373901296
$ws.Cells.Item(3, 6).Value = 33.99471573617401
$ws.Cells.Item(3, 7).Value = 35.59641863424952
$ws.Cells.Item(3, 8).Value = 16.19444284497173
$ws.Cells.Item(3, 9).Value = 26.12881362803258
$ws.Cells.Item(3, 10).Value = 10.97034501920068
$ws.Cells.Item(3, 11).Value = 8.164907029856039
$ws.Cells.Item(3, 12).Value = 11.66285455360478
$ws.Cells.Item(3, 13).Value = 14.59411688561036
$ws.Cells.Item(3, 15).Value = 25.51901394877637

$ws.Cells.Item(4, 2).Value = 10.92520777436292
$ws.Cells.Item(4, 3).Value = 9.892176296860681
$ws.Cells.Item(4, 4).Value = 9.571117873090158
$ws.Cells.Item(4, 6).Value = 34.06586760905979
$ws.Cells.Item(4, 7).Value = 35.69248845722553
$ws.Cells.Item(4, 8).Value = 16.22755935203281
$ws.Cells.Item(4, 9).Value = 26.19483628715333
$ws.Cells.Item(4, 10).Value = 10.98444201173618
$ws.Cells.Item(4, 11).Value = 8.041762038398742
$ws.Cells.Item(4, 12).Value = 11.66781448735228
$ws.Cells.Item(4, 13).Value = 14.57005952585583
$ws.Cells.Item(4, 15).Value = 25.57981719421065

$ws.Cells.Item(5, 2).Value = 10.8605026842988
$ws.Cells.Item(5, 3).Value = 9.890316714287302
$ws.Cells.Item(5, 4).Value = 9.569493567930889
$ws.Cells.Item(5, 6).Value = 34.09644087919796
$ws.Cells.Item(5, 7).Value = 35.73372328091618
$ws.Cells.Item(5, 8).Value = 16.24157373763847
$ws.Cells.Item(5, 9).Value = 26.2227348497478
$ws.Cells.Item(5, 10).Value = 10.99037222945737
$ws.Cells.Item(5, 11).Value = 7.991176813485596
$ws.Cells.Item(5, 12).Value = 11.67014707063976
$ws.Cells.Item(5, 13).Value = 14.5606273581853
$ws.Cells.Item(5, 15).Value = 25.60565534989157

$ws.Cells.Item(6, 2).Value = 10.84973832229361
$ws.Cells.Item(6, 3).Value = 9.890013064677694
$ws.Cells.Item(6, 4).Value = 9.569243195009056
$ws.Cells.Item(6, 6).Value = 34.10161284265658
$ws.Cells.Item(6, 7).Value = 35.74069613507901
$ws.Cells.Item(6, 8).Value = 16.24393219012627
$ws.Cells.Item(6, 9).Value = 26.22742743862392
$ws.Cells.Item(6, 10).Value = 10.99136816073579
$ws.Cells.Item(6, 11).Value = 7.982754847497656
$ws.Cells.Item(6, 12).Value = 11.67055322434995
$ws.Cells.Item(6, 13).Value = 14.55908375819112
$ws.Cells.Item(6, 15).Value = 25.61000980295849

$ws.Cells.Item(7, 2).Value = 10.9243365454944
$ws.Cells.Item(7, 3).Value = 9.892150874331387
$ws.Cells.Item(7, 4).Value = 9.571094671706563
$ws.Cells.Item(7, 6).Value = 34.06627354173774
$ws.Cells.Item(7, 7).Value = 35.69303612547319
$ws.Cells.Item(7, 8).Value = 16.22774625207231
$ws.Cells.Item(7, 9).Value = 26.19520851168737
$ws.Cells.Item(7, 10).Value = 10.98452123666148
$ws.Cells.Item(7, 11).Value = 8.041081369110696
$ws.Cells.Item(7, 12).Value = 11.66784468356255
$ws.Cells.Item(7, 13).Value = 14.5699308087118
$ws.Cells.Item(7, 15).Value = 25.58016136310912

$ws.Cells.Item(8, 2).Value = 11.24998717261715
$ws.Cells.Item(8, 3).Value = 9.902324719526847
$ws.Cells.Item(8, 4).Value = 9.58206265216389
$ws.Cells.Item(8, 6).Value = 33.92418784505877
$ws.Cells.Item(8, 7).Value = 35.50104371619253
$ws.Cells.Item(8, 8).Value = 16.16087709731183
$ws.Cells.Item(8, 9).Value = 26.06175532674203
$ws.Cells.Item(8, 10).Value = 10.95593515699898
$ws.Cells.Item(8, 11).Value = 8.294750944358452
$ws.Cells.Item(8, 12).Value = 11.6586237786041
$ws.Cells.Item(8, 13).Value = 14.62104080473797
$ws.Cells.Item(8, 15).Value = 25.45775520483048

$ws.Cells.Item(9, 2).Value = 11.87039068352641
$ws.Cells.Item(9, 3).Value = 9.92487990396501
$ws.Cells.Item(9, 4).Value = 9.613487421960528
$ws.Cells.Item(9, 6).Value = 33.69543364556018
$ws.Cells.Item(9, 7).Value = 35.19066442004505
$ws.Cells.Item(9, 8).Value = 16.04607698675585
$ws.Cells.Item(9, 9).Value = 25.83133184152172
$ws.Cells.Item(9, 10).Value = 10.90569533590507
$ws.Cells.Item(9, 11).Value = 8.774677625495181
$ws.Cells.Item(9, 12).Value = 11.65030109678377
$ws.Cells.Item(9, 13).Value = 14.73251712971418
$ws.Cells.Item(9, 15).Value = 25.25116237293947

$ws.Cells.Item(10, 2).Value = 12.31058206263785
$ws.Cells.Item(10, 3).Value = 9.942989760780724
$ws.Cells.Item(10, 4).Value = 9.642485950419131
$ws.Cells.Item(10, 6).Value = 33.557827073758
$ws.Cells.Item(10, 7).Value = 35.00317429638197
$ws.Cells.Item(10, 8).Value = 15.97164699467459
$ws.Cells.Item(10, 9).Value = 25.6810411666647
$ws.Cells.Item(10, 10).Value = 10.87229876700462
$ws.Cells.Item(10, 11).Value = 9.11319897482937
$ws.Cells.Item(10, 12).Value = 11.65008275720949
$ws.Cells.Item(10, 13).Value = 14.82089143141419
$ws.Cells.Item(10, 15).Value = 25.11975678627591

$ws.Cells.Item(11, 2).Value = 12.50657901795763
$ws.Cells.Item(11, 3).Value = 9.951551808233626
$ws.Cells.Item(11, 4).Value = 9.656929843106116
$ws.Cells.Item(11, 6).Value = 33.50185005425134
$ws.Cells.Item(11, 7).Value = 34.92673313044366
$ws.Cells.Item(11, 8).Value = 15.93993085597388
$ws.Cells.Item(11, 9).Value = 25.61678245482801
$ws.Cells.Item(11, 10).Value = 10.85786222024164
$ws.Cells.Item(11, 11).Value = 9.263500648499715
$ws.Cells.Item(11, 12).Value = 11.6512514029179
$ws.Cells.Item(11, 13).Value = 14.86241352469757
$ws.Cells.Item(11, 15).Value = 25.06439988978362

$ws.Cells.Item(12, 2).Value = 12.58012139541361
$ws.Cells.Item(12, 3).Value = 9.954839631059595
$ws.Cells.Item(12, 4).Value = 9.662576495430336
$ws.Cells.Item(12, 6).Value = 33.48160578986132
$ws.Cells.Item(12, 7).Value = 34.8990632164663
$ws.Cells.Item(12, 8).Value = 15.92822824021833
$ws.Cells.Item(12, 9).Value = 25.59303935511766
$ws.Cells.Item(12, 10).Value = 10.85250364071122
$ws.Cells.Item(12, 11).Value = 9.319836560963775
$ws.Cells.Item(12, 12).Value = 11.65187511616354
$ws.Cells.Item(12, 13).Value = 14.87831892049697
$ws.Cells.Item(12, 15).Value = 25.04407324848748

$ws.Cells.Item(13, 2).Value = 12.56431386337275
$ws.Cells.Item(13, 3).Value = 9.954129529302175
$ws.Cells.Item(13, 4).Value = 9.661352562486185
$ws.Cells.Item(13, 6).Value = 33.48592335491269
$ws.Cells.Item(13, 7).Value = 34.90496558065119
$ws.Cells.Item(13, 8).Value = 15.9307349349768
$ws.Cells.Item(13, 9).Value = 25.59812661325674
$ws.Cells.Item(13, 10).Value = 10.85365290027383
$ws.Cells.Item(13, 11).Value = 9.307730150846449
$ws.Cells.Item(13, 12).Value = 11.65173274974433
$ws.Cells.Item(13, 13).Value = 14.87488544871101
$ws.Cells.Item(13, 15).Value = 25.04842267310392

$ws.Cells.Item(14, 2).Value = 12.51264333345341
$ws.Cells.Item(14, 3).Value = 9.951821392053812
$ws.Cells.Item(14, 4).Value = 9.657390865779975
$ws.Cells.Item(14, 6).Value = 33.50016544154845
$ws.Cells.Item(14, 7).Value = 34.92443109947855
$ws.Cells.Item(14, 8).Value = 15.93896191273839
$ws.Cells.Item(14, 9).Value = 25.61481726911771
$ws.Cells.Item(14, 10).Value = 10.85741920016839
$ws.Cells.Item(14, 11).Value = 9.268147329869686
$ws.Cells.Item(14, 12).Value = 11.6512990920637
$ws.Cells.Item(14, 13).Value = 14.86371847236604
$ws.Cells.Item(14, 15).Value = 25.0627148609332

$ws.Cells.Item(15, 2).Value = 12.48090353270162
$ws.Cells.Item(15, 3).Value = 9.950413494868943
$ws.Cells.Item(15, 4).Value = 9.65498717898967
$ws.Cells.Item(15, 6).Value = 33.5090132679179
$ws.Cells.Item(15, 7).Value = 34.93652067784216
$ws.Cells.Item(15, 8).Value = 15.94404121834487
$ws.Cells.Item(15, 9).Value = 25.62511763154972
$ws.Cells.Item(15, 10).Value = 10.85974024890063
$ws.Cells.Item(15, 11).Value = 9.243824780871218
$ws.Cells.Item(15, 12).Value = 11.65105702281168
$ws.Cells.Item(15, 13).Value = 14.85690182327698
$ws.Cells.Item(15, 15).Value = 25.07155204649965

$ws.Cells.Item(16, 2).Value = 12.29768191513495
$ws.Cells.Item(16, 3).Value = 9.942436641208722
$ws.Cells.Item(16, 4).Value = 9.641566953268978
$ws.Cells.Item(16, 6).Value = 33.56161861603638
$ws.Cells.Item(16, 7).Value = 35.00834833550808
$ws.Cells.Item(16, 8).Value = 15.97376278875798
$ws.Cells.Item(16, 9).Value = 25.68532325917415
$ws.Cells.Item(16, 10).Value = 10.87325739823253
$ws.Cells.Item(16, 11).Value = 9.10329781095354
$ws.Cells.Item(16, 12).Value = 11.65003180040367
$ws.Cells.Item(16, 13).Value = 14.81820375475771
$ws.Cells.Item(16, 15).Value = 25.12346345151116

$ws.Cells.Item(17, 2).Value = 12.18414439907962
$ws.Cells.Item(17, 3).Value = 9.93762532317494
$ws.Cells.Item(17, 4).Value = 9.633652727643447
$ws.Cells.Item(17, 6).Value = 33.59558676923635
$ws.Cells.Item(17, 7).Value = 35.05468177565258
$ws.Cells.Item(17, 8).Value = 15.99254439236582
$ws.Cells.Item(17, 9).Value = 25.72330948687362
$ws.Cells.Item(17, 10).Value = 10.88174297207
$ws.Cells.Item(17, 11).Value = 9.01610701008843
$ws.Cells.Item(17, 12).Value = 11.64972680994607
$ws.Cells.Item(17, 13).Value = 14.79479627736466
$ws.Cells.Item(17, 15).Value = 25.15644166620653

$ws.Cells.Item(18, 2).Value = 12.1184444148765
$ws.Cells.Item(18, 3).Value = 9.934888516995183
$ws.Cells.Item(18, 4).Value = 9.629218738461214
$ws.Cells.Item(18, 6).Value = 33.61574749367971
$ws.Cells.Item(18, 7).Value = 35.08216421507787
$ws.Cells.Item(18, 8).Value = 16.00354875849798
$ws.Cells.Item(18, 9).Value = 25.74554497975896
$ws.Cells.Item(18, 10).Value = 10.88669480792204
$ws.Cells.Item(18, 11).Value = 8.965612203959452
$ws.Cells.Item(18, 12).Value = 11.6496708078176
$ws.Cells.Item(18, 13).Value = 14.78145766736109
$ws.Cells.Item(18, 15).Value = 25.17582589722325

$ws.Cells.Item(19, 2).Value = 12.09613358768997
$ws.Cells.Item(19, 3).Value = 9.933967156309299
$ws.Cells.Item(19, 4).Value = 9.62773783305161
$ws.Cells.Item(19, 6).Value = 33.62268056317478
$ws.Cells.Item(19, 7).Value = 35.09161217773649
$ws.Cells.Item(19, 8).Value = 16.00730930218928
$ws.Cells.Item(19, 9).Value = 25.7531399918848
$ws.Cells.Item(19, 10).Value = 10.88838364955678
$ws.Cells.Item(19, 11).Value = 8.948457838994049
$ws.Cells.Item(19, 12).Value = 11.64967239141531
$ws.Cells.Item(19, 13).Value = 14.77696310850781
$ws.Cells.Item(19, 15).Value = 25.18246051312446

$ws.Cells.Item(20, 2).Value = 12.19627216500979
$ws.Cells.Item(20, 3).Value = 9.938134342942679
$ws.Cells.Item(20, 4).Value = 9.634483012682852
$ws.Cells.Item(20, 6).Value = 33.59190630140314
$ws.Cells.Item(20, 7).Value = 35.04966329257713
$ws.Cells.Item(20, 8).Value = 15.9905241885409
$ws.Cells.Item(20, 9).Value = 25.71922575838396
$ws.Cells.Item(20, 10).Value = 10.88083230687315
$ws.Cells.Item(20, 11).Value = 9.025424689304968
$ws.Cells.Item(20, 12).Value = 11.64974692537045
$ws.Cells.Item(20, 13).Value = 14.7972751889327
$ws.Cells.Item(20, 15).Value = 25.15288801956042

$ws.Cells.Item(21, 2).Value = 12.5278391110563
$ws.Cells.Item(21, 3).Value = 9.952498119428704
$ws.Cells.Item(21, 4).Value = 9.658549731459628
$ws.Cells.Item(21, 6).Value = 33.49595632389968
$ws.Cells.Item(21, 7).Value = 34.91867892275216
$ws.Cells.Item(21, 8).Value = 15.93653710561963
$ws.Cells.Item(21, 9).Value = 25.60989880372098
$ws.Cells.Item(21, 10).Value = 10.85631001285038
$ws.Cells.Item(21, 11).Value = 9.279789877542909
$ws.Cells.Item(21, 12).Value = 11.65142156039685
$ws.Cells.Item(21, 13).Value = 14.86699361230314
$ws.Cells.Item(21, 15).Value = 25.0584996414725

$ws.Cells.Item(22, 2).Value = 12.74056140365699
$ws.Cells.Item(22, 3).Value = 9.962150769448765
$ws.Cells.Item(22, 4).Value = 9.675309283138393
$ws.Cells.Item(22, 6).Value = 33.43880274083715
$ws.Cells.Item(22, 7).Value = 34.8405164718116
$ws.Cells.Item(22, 8).Value = 15.90304612498974
$ws.Cells.Item(22, 9).Value = 25.54188765361788
$ws.Cells.Item(22, 10).Value = 10.8409139053347
$ws.Cells.Item(22, 11).Value = 9.442631187443855
$ws.Cells.Item(22, 12).Value = 11.65357146437324
$ws.Cells.Item(22, 13).Value = 14.91361493384658
$ws.Cells.Item(22, 15).Value = 25.00051748972949

$ws.Cells.Item(23, 2).Value = 12.62741200482818
$ws.Cells.Item(23, 3).Value = 9.9569750464418
$ws.Cells.Item(23, 4).Value = 9.666271133467905
$ws.Cells.Item(23, 6).Value = 33.46879809449519
$ws.Cells.Item(23, 7).Value = 34.88155081907949
$ws.Cells.Item(23, 8).Value = 15.92075701589911
$ws.Cells.Item(23, 9).Value = 25.57787189251289
$ws.Cells.Item(23, 10).Value = 10.84907354066638
$ws.Cells.Item(23, 11).Value = 9.35604607611784
$ws.Cells.Item(23, 12).Value = 11.6523278460095
$ws.Cells.Item(23, 13).Value = 14.88863828811185
$ws.Cells.Item(23, 15).Value = 25.03112447764913

$ws.Cells.Item(24, 2).Value = 12.19079052305799
$ws.Cells.Item(24, 3).Value = 9.937904123816578
$ws.Cells.Item(24, 4).Value = 9.634107279259243
$ws.Cells.Item(24, 6).Value = 33.59356827242474
$ws.Cells.Item(24, 7).Value = 35.05192951784501
$ws.Cells.Item(24, 8).Value = 15.99143687919657
$ws.Cells.Item(24, 9).Value = 25.72107077628305
$ws.Cells.Item(24, 10).Value = 10.88124379004678
$ws.Cells.Item(24, 11).Value = 9.021213307930385
$ws.Cells.Item(24, 12).Value = 11.64973745937068
$ws.Cells.Item(24, 13).Value = 14.79615410279952
$ws.Cells.Item(24, 15).Value = 25.15449330057833

$ws.Cells.Item(25, 2).Value = 11.70498801898848
$ws.Cells.Item(25, 3).Value = 9.918504324250042
$ws.Cells.Item(25, 4).Value = 9.603937836816906
$ws.Cells.Item(25, 6).Value = 33.75197230037421
$ws.Cells.Item(25, 7).Value = 35.26752682093635
$ws.Cells.Item(25, 8).Value = 16.07538955613201
$ws.Cells.Item(25, 9).Value = 25.89032597440346
$ws.Cells.Item(25, 10).Value = 10.91866705132892
$ws.Cells.Item(25, 11).Value = 8.647086134966191
$ws.Cells.Item(25, 12).Value = 11.65151307092506
$ws.Cells.Item(25, 13).Value = 14.70119423354038
$ws.Cells.Item(25, 15).Value = 25.30347184473472

Write-Host "Updated loading_percent values for 380 kV case"